# Fruta / hortaliza, semanal
# Insert 4 new weekly rows of grape (Uva) price data before row 813,
# pushing the existing rows 813-823 down to 817-827.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows starting at row 813 (existing data shifts down).
$ws.Rows("813:816").Insert()

# Fill the newly inserted rows with the new weekly price records.
$colVals = @{
    813 = @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44595, 13, "Fruta", 100109, "Uva", 100109001, "Uva", "Flame Seedless", "Primera", 1080, 8000, 9000, 8556, "`$/bandeja 18 kilos", "Región Metropolitana", 475, 18)
    814 = @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44595, 13, "Fruta", 100109, "Uva", 100109001, "Uva", "Red Globe", "Primera", 300, 12000, 13000, 12400, "`$/bandeja 18 kilos", "Provincia del Elquí", 689, 18)
    815 = @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44595, 13, "Fruta", 100109, "Uva", 100109001, "Uva", "Superior Seedless", "Primera", 400, 10000, 10000, 10000, "`$/bandeja 18 kilos", "Provincia del Elquí", 556, 18)
    816 = @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44595, 13, "Fruta", 100109, "Uva", 100109001, "Uva", "Superior Seedless", "Primera", 940, 10000, 11000, 10553, "`$/bandeja 18 kilos", "Región de O'Higgins", 586, 18)
}

foreach ($r in $colVals.Keys) {
    $vals = $colVals[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $vals[$i]
    }
}
